# edit.ps1 - Add an OBJECTIVE section to the resume, tidy the contact-info
# line and the "Web Engineering I and II" line, and relocate the trailing
# _GoBack bookmark into the new OBJECTIVE paragraph.

$d = $word.ActiveDocument

$bullet = [char]0x2022   # "•"

# Namespace declaration used for the InsertXML "mini packages" below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 0. Remove the old _GoBack bookmark that currently sits at the very end
#    of the document (after "Fluent in Spanish") BEFORE we create the
#    new one below - Bookmarks("_GoBack").Delete() removes the first
#    match by document order, so the old one must be gone before a new
#    "_GoBack" bookmark is introduced in the OBJECTIVE paragraph.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 1. Contact-info line: merge the split "github.com/<name>" + spell-check
#    wrapped runs + " linkedin.com/..." run back into a single plain run.
# ---------------------------------------------------------------------
$contactPara = $d.Paragraphs(3)
$contactText = "(803)389-6750 $bullet danielmartincraig@gmail.com $bullet github.com/danielmartincraig $bullet linkedin.com/danielcraig23"
$contactBody = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/></w:rPr><w:t>' + $contactText + '</w:t></w:r></w:p>'
$contactPara.Range.InsertXML((New-PkgXml $contactBody))

# ---------------------------------------------------------------------
# 2. Insert a brand-new "OBJECTIVE:" paragraph (Heading1 style) right
#    before the "EDUCATION:" heading. Splitting the EDUCATION paragraph
#    keeps the new paragraph free of any stray direct formatting, which
#    matches how the target document renders the "OBJECTIVE: " run.
# ---------------------------------------------------------------------
$eduPara = $d.Paragraphs(4)
$eduPara.Range.InsertParagraphBefore()

$objectivePara = $d.Paragraphs(4)
$objectiveBody = '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Eager to drive back-end solutions at </w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Pariveda</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> on a full-time basis</w:t></w:r>' + `
    '</w:p>'
$objectivePara.Range.InsertXML((New-PkgXml $objectiveBody))

# ---------------------------------------------------------------------
# 3. "Web Engineering I and II" bullet: collapse the gramStart/gramEnd
#    proof-error wrapped runs back into one plain run.
# ---------------------------------------------------------------------
$webPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Web*Engineering I and II*") {
        $webPara = $p
        break
    }
}
$webBody = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>' + $bullet + '    Web Engineering I and II</w:t></w:r></w:p>'
$webPara.Range.InsertXML((New-PkgXml $webBody))

Write-Output "edit complete"
